# Generate Report for Handoff
# Updates the localization-status report to reflect that the files are now
# "Ready for handoff" instead of "In Translation", and refreshes the
# handoff timestamps accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (row 2 = the .md file) ---
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-30-18 22:30:25"

# --- zh-cn sheet (row 2) ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-18 22:30:23"

# --- de-de sheet (row 2) ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-18 22:30:25"

$wb.Save()
